# Updates per-row market/profit figures (columns H-N: currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 114.666664
$ws.Cells.Item(2, 9).Value = 114.5
$ws.Cells.Item(2, 10).Value = 115
$ws.Cells.Item(2, 11).Value = 114.5
$ws.Cells.Item(2, 12).Value = 115
$ws.Cells.Item(2, 13).Value = -1.5
$ws.Cells.Item(2, 14).Value = -341

# Row 12
$ws.Cells.Item(12, 8).Value = 286.83334
$ws.Cells.Item(12, 9).Value = 192.75
$ws.Cells.Item(12, 11).Value = 192.75
$ws.Cells.Item(12, 13).Value = -22.75

# Row 15
$ws.Cells.Item(15, 8).Value = 98.69
$ws.Cells.Item(15, 9).Value = 98.69
$ws.Cells.Item(15, 11).Value = 296.07
$ws.Cells.Item(15, 13).Value = -127.07

# Row 21
$ws.Cells.Item(21, 8).Value = 8346.75
$ws.Cells.Item(21, 9).Value = 7368
$ws.Cells.Item(21, 10).Value = 8673
$ws.Cells.Item(21, 11).Value = 7368
$ws.Cells.Item(21, 12).Value = 8673
$ws.Cells.Item(21, 13).Value = -6900
$ws.Cells.Item(21, 14).Value = -9609

# Row 23
$ws.Cells.Item(23, 8).Value = 8346.75
$ws.Cells.Item(23, 9).Value = 7368
$ws.Cells.Item(23, 10).Value = 8673
$ws.Cells.Item(23, 11).Value = 7368
$ws.Cells.Item(23, 12).Value = 8673
$ws.Cells.Item(23, 13).Value = -7134
$ws.Cells.Item(23, 14).Value = -9141

# Row 38
$ws.Cells.Item(38, 8).Value = 797.53845
$ws.Cells.Item(38, 9).Value = 99
$ws.Cells.Item(38, 10).Value = 1612.5
$ws.Cells.Item(38, 11).Value = 297
$ws.Cells.Item(38, 12).Value = 4837.5
$ws.Cells.Item(38, 13).Value = 75
$ws.Cells.Item(38, 14).Value = -5581.5

# Row 58
$ws.Cells.Item(58, 8).Value = 1593.7
$ws.Cells.Item(58, 9).Value = 248
$ws.Cells.Item(58, 10).Value = 2939.4
$ws.Cells.Item(58, 11).Value = 744
$ws.Cells.Item(58, 12).Value = 8818.200000000001
$ws.Cells.Item(58, 13).Value = -594
$ws.Cells.Item(58, 14).Value = -9118.200000000001

# Row 87
$ws.Cells.Item(87, 8).Value = 26968
$ws.Cells.Item(87, 10).Value = 26968
$ws.Cells.Item(87, 12).Value = 26968
$ws.Cells.Item(87, 14).Value = -29464

# Row 90
$ws.Cells.Item(90, 8).Value = 26968
$ws.Cells.Item(90, 10).Value = 26968
$ws.Cells.Item(90, 12).Value = 80904
$ws.Cells.Item(90, 14).Value = -93384

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 7211.8315
$ws.Cells.Item(32, 9).Value = 4419.122
$ws.Cells.Item(32, 11).Value = 4419.122
$ws.Cells.Item(32, 13).Value = -4132.122

# Row 61
$ws.Cells.Item(61, 8).Value = 3549.9429
$ws.Cells.Item(61, 9).Value = 2415.111
$ws.Cells.Item(61, 10).Value = 4751.5293
$ws.Cells.Item(61, 11).Value = 2415.111
$ws.Cells.Item(61, 12).Value = 4751.5293
$ws.Cells.Item(61, 13).Value = -2203.111
$ws.Cells.Item(61, 14).Value = -5175.5293

# Row 133
$ws.Cells.Item(133, 8).Value = 22565.25
$ws.Cells.Item(133, 10).Value = 22565.25
$ws.Cells.Item(133, 12).Value = 22565.25
$ws.Cells.Item(133, 14).Value = -27625.25

# Row 136
$ws.Cells.Item(136, 8).Value = 3549.9429
$ws.Cells.Item(136, 9).Value = 2415.111
$ws.Cells.Item(136, 10).Value = 4751.5293
$ws.Cells.Item(136, 11).Value = 7245.333
$ws.Cells.Item(136, 12).Value = 14254.5879
$ws.Cells.Item(136, 13).Value = -4695.333
$ws.Cells.Item(136, 14).Value = -19354.5879

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1859.3636
$ws.Cells.Item(20, 9).Value = 1853
$ws.Cells.Item(20, 10).Value = 1867
$ws.Cells.Item(20, 11).Value = 1853
$ws.Cells.Item(20, 12).Value = 1867
$ws.Cells.Item(20, 13).Value = -1606
$ws.Cells.Item(20, 14).Value = -2361

# Row 26
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(26, 14).ClearContents()

# Row 29
$ws.Cells.Item(29, 8).Value = 2945
$ws.Cells.Item(29, 9).Value = 593.3333
$ws.Cells.Item(29, 10).Value = 10000
$ws.Cells.Item(29, 11).Value = 593.3333
$ws.Cells.Item(29, 12).Value = 10000
$ws.Cells.Item(29, 13).Value = -304.3333
$ws.Cells.Item(29, 14).Value = -10578

# Row 86
$ws.Cells.Item(86, 8).Value = 7858.4707
$ws.Cells.Item(86, 9).Value = 8350
$ws.Cells.Item(86, 11).Value = 8350
$ws.Cells.Item(86, 13).Value = -7227

# Row 89
$ws.Cells.Item(89, 8).Value = 7858.4707
$ws.Cells.Item(89, 9).Value = 8350
$ws.Cells.Item(89, 11).Value = 41750
$ws.Cells.Item(89, 13).Value = -36134

# Row 134
$ws.Cells.Item(134, 8).Value = 34969
$ws.Cells.Item(134, 9).Value = 47412.87
$ws.Cells.Item(134, 10).Value = 8950
$ws.Cells.Item(134, 11).Value = 142238.61
$ws.Cells.Item(134, 12).Value = 26850
$ws.Cells.Item(134, 13).Value = -139703.61
$ws.Cells.Item(134, 14).Value = -31920

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 3322.6885
$ws.Cells.Item(31, 9).Value = 2730.8696
$ws.Cells.Item(31, 10).Value = 3680.8948
$ws.Cells.Item(31, 11).Value = 2730.8696
$ws.Cells.Item(31, 12).Value = 3680.8948
$ws.Cells.Item(31, 13).Value = -2435.8696
$ws.Cells.Item(31, 14).Value = -4270.8948

# Row 33
$ws.Cells.Item(33, 8).Value = 1050
$ws.Cells.Item(33, 9).Value = 1050
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 1050
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = -671
$ws.Cells.Item(33, 14).ClearContents()

# Row 34
$ws.Cells.Item(34, 8).Value = 3322.6885
$ws.Cells.Item(34, 9).Value = 2730.8696
$ws.Cells.Item(34, 10).Value = 3680.8948
$ws.Cells.Item(34, 11).Value = 2730.8696
$ws.Cells.Item(34, 12).Value = 3680.8948
$ws.Cells.Item(34, 13).Value = -2528.8696
$ws.Cells.Item(34, 14).Value = -4084.8948

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Cells.Item(113, 8).Value = 525.4942600000001
$ws.Cells.Item(113, 9).Value = 530.2292
$ws.Cells.Item(113, 10).Value = 519.6667
$ws.Cells.Item(113, 11).Value = 1590.6876
$ws.Cells.Item(113, 12).Value = 1559.0001
$ws.Cells.Item(113, 13).Value = 579.3124
$ws.Cells.Item(113, 14).Value = -5899.0001

# Row 131
$ws.Cells.Item(131, 8).Value = 867.1
$ws.Cells.Item(131, 9).Value = 692
$ws.Cells.Item(131, 10).Value = 874.3958
$ws.Cells.Item(131, 11).Value = 2076
$ws.Cells.Item(131, 12).Value = 2623.1874
$ws.Cells.Item(131, 13).Value = 2964
$ws.Cells.Item(131, 14).Value = -12703.1874

# Row 134
$ws.Cells.Item(134, 8).Value = 2728.6538
$ws.Cells.Item(134, 9).Value = 2070.6155
$ws.Cells.Item(134, 10).Value = 3386.6924
$ws.Cells.Item(134, 11).Value = 6211.8465
$ws.Cells.Item(134, 12).Value = 10160.0772
$ws.Cells.Item(134, 13).Value = -1141.8465
$ws.Cells.Item(134, 14).Value = -20300.0772

# Row 139
$ws.Cells.Item(139, 8).Value = 83334210
$ws.Cells.Item(139, 9).Value = 100000650
$ws.Cells.Item(139, 10).Value = 2000
$ws.Cells.Item(139, 11).Value = 300001950
$ws.Cells.Item(139, 12).Value = 6000
$ws.Cells.Item(139, 13).Value = -299996810
$ws.Cells.Item(139, 14).Value = -16280

# Row 140
$ws.Cells.Item(140, 8).Value = 1555.1305
$ws.Cells.Item(140, 9).Value = 1169.8182
$ws.Cells.Item(140, 11).Value = 3509.4546
$ws.Cells.Item(140, 13).Value = 1670.5454

$ws = $wb.Worksheets.Item("GSM")
# Row 27
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 13).ClearContents()

# Row 29
$ws.Cells.Item(29, 8).Value = 8666.666999999999
$ws.Cells.Item(29, 10).Value = 10000
$ws.Cells.Item(29, 12).Value = 10000
$ws.Cells.Item(29, 14).Value = -10580

# Row 107
$ws.Cells.Item(107, 8).Value = 3628.3667
$ws.Cells.Item(107, 9).Value = 5704.8335
$ws.Cells.Item(107, 10).Value = 513.6667
$ws.Cells.Item(107, 11).Value = 5704.8335
$ws.Cells.Item(107, 12).Value = 513.6667
$ws.Cells.Item(107, 13).Value = -3784.8335
$ws.Cells.Item(107, 14).Value = -4353.6667

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 321
$ws.Cells.Item(22, 9).Value = 235.6
$ws.Cells.Item(22, 10).Value = 748
$ws.Cells.Item(22, 11).Value = 235.6
$ws.Cells.Item(22, 12).Value = 748
$ws.Cells.Item(22, 13).Value = 59.40000000000001
$ws.Cells.Item(22, 14).Value = -1338

# Row 27
$ws.Cells.Item(27, 8).Value = 321
$ws.Cells.Item(27, 9).Value = 235.6
$ws.Cells.Item(27, 10).Value = 748
$ws.Cells.Item(27, 11).Value = 235.6
$ws.Cells.Item(27, 12).Value = 748
$ws.Cells.Item(27, 13).Value = -128.6
$ws.Cells.Item(27, 14).Value = -962

$ws = $wb.Worksheets.Item("WVR")
# Row 33
$ws.Cells.Item(33, 8).Value = 3716.6667
$ws.Cells.Item(33, 9).Value = 3500
$ws.Cells.Item(33, 10).Value = 3760
$ws.Cells.Item(33, 11).Value = 3500
$ws.Cells.Item(33, 12).Value = 3760
$ws.Cells.Item(33, 13).Value = -3250
$ws.Cells.Item(33, 14).Value = -4260

# Row 36
$ws.Cells.Item(36, 8).Value = 3716.6667
$ws.Cells.Item(36, 9).Value = 3500
$ws.Cells.Item(36, 10).Value = 3760
$ws.Cells.Item(36, 11).Value = 3500
$ws.Cells.Item(36, 12).Value = 3760
$ws.Cells.Item(36, 13).Value = -3250
$ws.Cells.Item(36, 14).Value = -4260

# Row 54
$ws.Cells.Item(54, 8).Value = 9317.5
$ws.Cells.Item(54, 10).Value = 9317.5
$ws.Cells.Item(54, 12).Value = 9317.5
$ws.Cells.Item(54, 14).Value = -10357.5

# Row 81
$ws.Cells.Item(81, 8).Value = 4435
$ws.Cells.Item(81, 9).Value = 3390
$ws.Cells.Item(81, 10).Value = 4783.3335
$ws.Cells.Item(81, 11).Value = 6780
$ws.Cells.Item(81, 12).Value = 9566.666999999999
$ws.Cells.Item(81, 13).Value = -5719
$ws.Cells.Item(81, 14).Value = -11688.667

# Row 84
$ws.Cells.Item(84, 8).Value = 4435
$ws.Cells.Item(84, 9).Value = 3390
$ws.Cells.Item(84, 10).Value = 4783.3335
$ws.Cells.Item(84, 11).Value = 33900
$ws.Cells.Item(84, 12).Value = 47833.335
$ws.Cells.Item(84, 13).Value = -28596
$ws.Cells.Item(84, 14).Value = -58441.335

# Row 107
$ws.Cells.Item(107, 8).Value = 316.27777
$ws.Cells.Item(107, 9).Value = 307.33334
$ws.Cells.Item(107, 10).Value = 361
$ws.Cells.Item(107, 11).Value = 922.0000200000001
$ws.Cells.Item(107, 12).Value = 1083
$ws.Cells.Item(107, 13).Value = 997.9999799999999
$ws.Cells.Item(107, 14).Value = -4923
